$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 content swap: jsmith2024@rediffmail.com/redcow@1999 -> chowdhurygs@rediffmail.com/Shajlee1@ ---
# Set B4 first, then A4, so the new shared-string entries land in the same order as the target file.
$ws.Range("B4").Value = "Shajlee1@"
$ws.Range("A4").Value = "chowdhurygs@rediffmail.com"

# --- Bigger font across the whole used range (username/password table) ---
$used = $ws.UsedRange
$used.Font.Size = 36

# --- Row heights / column widths grow along with the font ---
$used.EntireRow.RowHeight = 46.5
$used.EntireColumn.ColumnWidth = 115.14

# --- Selection moves to A10 ---
$ws.Range("A10").Select()
